# Feria Lagunitas de Puerto Montt - Palta
# Weekly update: insert 3 new "Quillota" price rows (Especial/Primera/Segunda)
# ahead of the existing data, pushing all subsequent rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 637..639 (existing rows 637+ shift down to 640+)
$ws.Rows.Item(637).Insert()
$ws.Rows.Item(638).Insert()
$ws.Rows.Item(639).Insert()

# Common columns shared by every "Feria Lagunitas de Puerto Montt" / Palta / Hass row
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$variedad  = "Hass"
$fecha     = 44939

# Row 637 - Especial
$ws.Range("A637").Value = $mercadoId
$ws.Range("B637").Value = $mercado
$ws.Range("C637").Value = $region
$ws.Range("D637").Value = $fecha
$ws.Range("E637").Value = $codreg
$ws.Range("F637").Value = $tipo
$ws.Range("G637").Value = $productoId
$ws.Range("H637").Value = $producto
$ws.Range("I637").Value = $categoriaId
$ws.Range("J637").Value = $categoria
$ws.Range("K637").Value = $variedad
$ws.Range("L637").Value = "Especial"
$ws.Range("M637").Value = 200
$ws.Range("N637").Value = 4800
$ws.Range("O637").Value = 4800
$ws.Range("P637").Value = 4800
$ws.Range("Q637").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R637").Value = "Provincia de Quillota"
$ws.Range("S637").Value = 4800
$ws.Range("T637").Value = 1

# Row 638 - Primera
$ws.Range("A638").Value = $mercadoId
$ws.Range("B638").Value = $mercado
$ws.Range("C638").Value = $region
$ws.Range("D638").Value = $fecha
$ws.Range("E638").Value = $codreg
$ws.Range("F638").Value = $tipo
$ws.Range("G638").Value = $productoId
$ws.Range("H638").Value = $producto
$ws.Range("I638").Value = $categoriaId
$ws.Range("J638").Value = $categoria
$ws.Range("K638").Value = $variedad
$ws.Range("L638").Value = "Primera"
$ws.Range("M638").Value = 150
$ws.Range("N638").Value = 4500
$ws.Range("O638").Value = 4500
$ws.Range("P638").Value = 4500
$ws.Range("Q638").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R638").Value = "Provincia de Quillota"
$ws.Range("S638").Value = 4500
$ws.Range("T638").Value = 1

# Row 639 - Segunda
$ws.Range("A639").Value = $mercadoId
$ws.Range("B639").Value = $mercado
$ws.Range("C639").Value = $region
$ws.Range("D639").Value = $fecha
$ws.Range("E639").Value = $codreg
$ws.Range("F639").Value = $tipo
$ws.Range("G639").Value = $productoId
$ws.Range("H639").Value = $producto
$ws.Range("I639").Value = $categoriaId
$ws.Range("J639").Value = $categoria
$ws.Range("K639").Value = $variedad
$ws.Range("L639").Value = "Segunda"
$ws.Range("M639").Value = 150
$ws.Range("N639").Value = 4000
$ws.Range("O639").Value = 4000
$ws.Range("P639").Value = 4000
$ws.Range("Q639").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R639").Value = "Provincia de Quillota"
$ws.Range("S639").Value = 4000
$ws.Range("T639").Value = 1
